$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update 24VDC cable quantity (row 22, "CABLE-004/BYNAT"): Qty for 1 cable
# 50 -> 0.3. The dependent "Qty for Tester" formula in D22 ($C$16*C22)
# recalculates automatically (400 -> 2.4).
$ws.Range("C22").Value = 0.3

# Update the last-selected cell saved with the sheet view.
$ws.Range("F9").Select() | Out-Null
